# Generate Report for Handoff
# Updates status text + timestamps on all sheets, and narrows the
# "Latest Handoff Datetime" / duplicate status columns on each sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-01 19:09:47"
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-01 19:09:42"
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-01 19:09:47"
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
